# B6-PowerPoint.pptx edit
#
# 1. Re-theme the deck from the "Integral" (Red Violet) design to the
#    built-in "Office Theme" colour scheme (the slide master's theme,
#    serialised as ppt/theme/theme1.xml, picks up the classic Office
#    palette: dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
# 2. Re-point the three data tables (on slides 14, 15 and 16) at the new
#    table style GUID that goes with that theme.

$p = $ppt.ActivePresentation

function Set-ThemeRgb {
    param($ThemeColorScheme, [int]$Index, [string]$Hex)
    $r = [Convert]::ToInt32($Hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4,2), 16)
    $ThemeColorScheme.Colors($Index).RGB = $r + ($g * 256) + ($b * 65536)
}

# The presentation has a single slide master/theme shared by every slide,
# so updating it once via slide 1 re-colours the whole deck.
$tcs = $p.Slides.Item(1).ThemeColorScheme

Set-ThemeRgb $tcs 1  "000000"   # dk1
Set-ThemeRgb $tcs 2  "FFFFFF"   # lt1
Set-ThemeRgb $tcs 3  "44546A"   # dk2
Set-ThemeRgb $tcs 4  "E7E6E6"   # lt2
Set-ThemeRgb $tcs 5  "5B9BD5"   # accent1
Set-ThemeRgb $tcs 6  "ED7D31"   # accent2
Set-ThemeRgb $tcs 7  "A5A5A5"   # accent3
Set-ThemeRgb $tcs 8  "FFC000"   # accent4
Set-ThemeRgb $tcs 9  "4472C4"   # accent5
Set-ThemeRgb $tcs 10 "70AD47"   # accent6
Set-ThemeRgb $tcs 11 "0563C1"   # hlink
Set-ThemeRgb $tcs 12 "954F72"   # folHlink

# Re-apply the matching table style (Medium-style GUID for the new theme)
# to the one table that lives on each of these three slides.
$newTableStyle = "{C33D7077-AE20-41C4-B7B2-E8098E3F319E}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}
